$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13/14 and 38/39/40 swap Coin name + Link (B, C) along with new D/E values.
# All other rows keep B/C unchanged, only D (Price) and/or E (Volume(1h)) change.

$updates = @(
    @{ Row = 2;  D = "25.570.66";   E = "  +5.45%  " },
    @{ Row = 3;  D = "1.756.75";    E = "  +4.89%  " },
    @{ Row = 4;  D = "1.004";       E = "  +0.52%  " },
    @{ Row = 5;  D = "317.22";      E = "  +2.76%  " },
    @{ Row = 6;  D = "1.002";       E = "  +0.53%  " },
    @{ Row = 7;  D = "0.3823";      E = "  +2.37%  " },
    @{ Row = 8;  D = "0.3588";      E = "  +4.33%  " },
    @{ Row = 9;  D = "49.59";       E = "  +3.40%  " },
    @{ Row = 10; D = "1.226";       E = "  +3.16%  " },
    @{ Row = 11; D = "0.07676";     E = "  +5.17%  " },
    @{ Row = 12; D = "1.004";       E = "  +0.75%  " },
    @{ Row = 13; B = "Solana";   C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol";   D = "21.55";  E = "  +4.70%  " },
    @{ Row = 14; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "6.490";  E = "  +6.81%  " },
    @{ Row = 15; E = "  +5.35%  " },
    @{ Row = 16; D = "1.755.60";    E = "  +5.11%  " },
    @{ Row = 17; D = "0.00001156"; E = "  +4.10%  " },
    @{ Row = 18; D = "0.06780";     E = "  +0.88%  " },
    @{ Row = 19; D = "1.002";       E = "  +0.49%  " },
    @{ Row = 20; D = "86.31";       E = "  +5.17%  " },
    @{ Row = 21; D = "17.65";       E = "  +7.02%  " },
    @{ Row = 22; D = "6.535";       E = "  +5.94%  " },
    @{ Row = 23; D = "13.18";       E = "  +9.55%  " },
    @{ Row = 24; D = "25.536.01";   E = "  +5.75%  " },
    @{ Row = 25; D = "2.449";       E = "  +1.80%  " },
    @{ Row = 26; D = "2.903";       E = "  +8.68%  " },
    @{ Row = 27; D = "20.98";       E = "  +7.47%  " },
    @{ Row = 28; D = "154.58";      E = "  +1.70%  " },
    @{ Row = 29; D = "1.948.59";    E = "  +5.12%  " },
    @{ Row = 30; D = "134.14";      E = "  +5.11%  " },
    @{ Row = 31; D = "1.210";       E = "  +22.55%  " },
    @{ Row = 32; D = "7.189";       E = "  +13.66%  " },
    @{ Row = 33; D = "4.218";       E = "  +4.71%  " },
    @{ Row = 34; D = "14.38";       E = "  +16.21%  " },
    @{ Row = 35; D = "1.810";       E = "  +3.69%  " },
    @{ Row = 36; D = "0.08754";     E = "  +3.85%  " },
    @{ Row = 37; D = "5.784";       E = "  +7.75%  " },
    @{ Row = 38; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "9.405";   E = "  +5.10%  " },
    @{ Row = 39; B = "VeChain";   C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.02498"; E = "  +6.09%  " },
    @{ Row = 40; B = "Hedera";    C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.06719"; E = "  +4.84%  " },
    @{ Row = 41; D = "0.2264";     E = "  +6.82%  " },
    @{ Row = 42; D = "1.290";      E = "  -0.51%  " },
    @{ Row = 43; D = "0.6597";     E = "  +7.43%  " },
    @{ Row = 44; D = "14.37";      E = "  +9.45%  " },
    @{ Row = 45; E = "  +0.56%  " },
    @{ Row = 46; D = "0.6308";     E = "  +5.40%  " },
    @{ Row = 47; D = "3.904";      E = "  +2.49%  " },
    @{ Row = 48; D = "2.188";      E = "  +7.91%  " },
    @{ Row = 49; D = "131.79";     E = "  +3.68%  " },
    @{ Row = 50; D = "0.07456";    E = "  +4.68%  " },
    @{ Row = 51; D = "81.16";      E = "  +6.57%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Price column holds text-formatted numbers (e.g. "25.570.66", "1.004").
        # Force the cell to Text format first so Excel doesn't silently
        # reinterpret the string as a numeric value.
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
